$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance Sheet")

# Mark these participants absent ("A") for the 9th Feb session column(s).
# Setting WrapText explicitly nudges the cell onto the same "applyAlignment"
# style variant already used by the other "A" cells in this sheet.
$cells = @("P12", "P13", "O17", "O21", "O24", "O26", "O33", "O34", "P34", "P39", "P44", "P50", "P54", "O61", "O63", "O68", "P68")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value = "A"
    $rng.WrapText = $false
}
